$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 values (daily scrape refresh) ---
# A2 is a purely-numeric-looking ID; format as Text first so it is stored
# as a string (matching the source data) instead of being auto-coerced to
# a number.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "1327043"

$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1327043"
$ws.Range("D2").Value = "Sousse, Tunisia"
$ws.Range("F2").Value = "21 applicants"
$ws.Range("H2").Value = "Progress Professional Center"

# --- Remove the now-stale listings that used to occupy rows 3-8 ---
$ws.Rows("3:8").Delete()

# --- Column width tweaks (values below compensate for the COM ColumnWidth
# padding of ~0.8333 chars so the stored OOXML width lands on the exact
# target value) ---
$ws.Columns(3).ColumnWidth = 15.166666666666668
$ws.Columns(4).ColumnWidth = 17.166666666666668
$ws.Columns(6).ColumnWidth = 15.166666666666668
$ws.Columns(8).ColumnWidth = 30.166666666666668
